# Fixes config-excel: make sheet names consistent with model terminology
# (e.g. gridNode instead of netNode), bump a capacity value, and refresh
# the saved selection/window state.

$wb = $excel.ActiveWorkbook

# --- Rename worksheets so "net" -> "grid" ---------------------------------
$wsNodes = $wb.Worksheets.Item("config_netNodes")
$wsNodes.Name = "config_gridNodes"

$wsConnections = $wb.Worksheets.Item("config_netConnections")
$wsConnections.Name = "config_gridConnections"

# --- config_gridNodes: refresh the saved selection ------------------------
$wsNodes.Activate()
$wsNodes.Range("E35").Select()

# --- config_gridConnections: bump a capacity value and refresh selection --
$wsConnections.Activate()
$wsConnections.Range("H19").Value = 12
$wsConnections.Range("Q27").Select()

# --- Update the saved window position/size of the workbook ---------------
$excel.Left = -120
$excel.Top = -120
$excel.Width = 2904
$excel.Height = 1584
